$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D ("last charge end time") for existing data rows (2-50) to the new timestamp
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 4).Value2 = 45956.999976851854
}

# Update rows 20-52: station name (A), terminal name (B), last-not-charged time (C), last charge end time (D)
$ws.Cells.Item(20, 1).Value2 = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(20, 2).Value2 = "101号直流"
$ws.Cells.Item(20, 3).Value2 = 45954.028229166666
$ws.Cells.Item(20, 4).Value2 = 45956.999976851854

$ws.Cells.Item(21, 1).Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(21, 2).Value2 = "004A号直流"
$ws.Cells.Item(21, 3).Value2 = 45954.540092592593
$ws.Cells.Item(21, 4).Value2 = 45956.999976851854

$ws.Cells.Item(22, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(22, 2).Value2 = "602号直流"
$ws.Cells.Item(22, 3).Value2 = 45955.051435185182
$ws.Cells.Item(22, 4).Value2 = 45956.999976851854

$ws.Cells.Item(23, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(23, 2).Value2 = "403号直流"
$ws.Cells.Item(23, 3).Value2 = 45955.164687500001
$ws.Cells.Item(23, 4).Value2 = 45956.999976851854

$ws.Cells.Item(24, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(24, 2).Value2 = "801号直流"
$ws.Cells.Item(24, 3).Value2 = 45955.221736111111
$ws.Cells.Item(24, 4).Value2 = 45956.999976851854

$ws.Cells.Item(25, 1).Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(25, 2).Value2 = "108号直流"
$ws.Cells.Item(25, 3).Value2 = 45955.55228009259
$ws.Cells.Item(25, 4).Value2 = 45956.999976851854

$ws.Cells.Item(26, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(26, 2).Value2 = "502号直流"
$ws.Cells.Item(26, 3).Value2 = 45955.558877314812
$ws.Cells.Item(26, 4).Value2 = 45956.999976851854

$ws.Cells.Item(27, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(27, 2).Value2 = "802号直流"
$ws.Cells.Item(27, 3).Value2 = 45955.568888888891
$ws.Cells.Item(27, 4).Value2 = 45956.999976851854

$ws.Cells.Item(28, 1).Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(28, 2).Value2 = "905号直流"
$ws.Cells.Item(28, 3).Value2 = 45956.041655092595
$ws.Cells.Item(28, 4).Value2 = 45956.999976851854

$ws.Cells.Item(29, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(29, 2).Value2 = "601号直流"
$ws.Cells.Item(29, 3).Value2 = 45956.04414351852
$ws.Cells.Item(29, 4).Value2 = 45956.999976851854

$ws.Cells.Item(30, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(30, 2).Value2 = "B02号直流"
$ws.Cells.Item(30, 3).Value2 = 45956.071863425925
$ws.Cells.Item(30, 4).Value2 = 45956.999976851854

$ws.Cells.Item(31, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(31, 2).Value2 = "B04号直流"
$ws.Cells.Item(31, 3).Value2 = 45956.156759259262
$ws.Cells.Item(31, 4).Value2 = 45956.999976851854

$ws.Cells.Item(32, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(32, 2).Value2 = "205号直流"
$ws.Cells.Item(32, 3).Value2 = 45956.1799537037
$ws.Cells.Item(32, 4).Value2 = 45956.999976851854

$ws.Cells.Item(33, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(33, 2).Value2 = "703号直流"
$ws.Cells.Item(33, 3).Value2 = 45956.187754629631
$ws.Cells.Item(33, 4).Value2 = 45956.999976851854

$ws.Cells.Item(34, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(34, 2).Value2 = "305号直流"
$ws.Cells.Item(34, 3).Value2 = 45956.194756944446
$ws.Cells.Item(34, 4).Value2 = 45956.999976851854

$ws.Cells.Item(35, 1).Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(35, 2).Value2 = "402号直流"
$ws.Cells.Item(35, 3).Value2 = 45956.213599537034
$ws.Cells.Item(35, 4).Value2 = 45956.999976851854

$ws.Cells.Item(36, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(36, 2).Value2 = "B01号直流"
$ws.Cells.Item(36, 3).Value2 = 45956.221458333333
$ws.Cells.Item(36, 4).Value2 = 45956.999976851854

$ws.Cells.Item(37, 1).Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(37, 2).Value2 = "903号直流"
$ws.Cells.Item(37, 3).Value2 = 45956.249155092592
$ws.Cells.Item(37, 4).Value2 = 45956.999976851854

$ws.Cells.Item(38, 1).Value2 = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(38, 2).Value2 = "106号直流"
$ws.Cells.Item(38, 3).Value2 = 45956.24931712963
$ws.Cells.Item(38, 4).Value2 = 45956.999976851854

$ws.Cells.Item(39, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(39, 2).Value2 = "401号直流"
$ws.Cells.Item(39, 3).Value2 = 45956.251469907409
$ws.Cells.Item(39, 4).Value2 = 45956.999976851854

$ws.Cells.Item(40, 1).Value2 = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(40, 2).Value2 = "103号直流"
$ws.Cells.Item(40, 3).Value2 = 45956.258935185186
$ws.Cells.Item(40, 4).Value2 = 45956.999976851854

$ws.Cells.Item(41, 1).Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(41, 2).Value2 = "111号直流"
$ws.Cells.Item(41, 3).Value2 = 45956.342476851853
$ws.Cells.Item(41, 4).Value2 = 45956.999976851854

$ws.Cells.Item(42, 1).Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(42, 2).Value2 = "101号直流"
$ws.Cells.Item(42, 3).Value2 = 45956.344872685186
$ws.Cells.Item(42, 4).Value2 = 45956.999976851854

$ws.Cells.Item(43, 1).Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(43, 2).Value2 = "109号直流"
$ws.Cells.Item(43, 3).Value2 = 45956.350162037037
$ws.Cells.Item(43, 4).Value2 = 45956.999976851854

$ws.Cells.Item(44, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(44, 2).Value2 = "A04号直流"
$ws.Cells.Item(44, 3).Value2 = 45956.380520833336
$ws.Cells.Item(44, 4).Value2 = 45956.999976851854

$ws.Cells.Item(45, 1).Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(45, 2).Value2 = "501号直流"
$ws.Cells.Item(45, 3).Value2 = 45956.392893518518
$ws.Cells.Item(45, 4).Value2 = 45956.999976851854

$ws.Cells.Item(46, 1).Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(46, 2).Value2 = "006A号直流"
$ws.Cells.Item(46, 3).Value2 = 45956.404502314814
$ws.Cells.Item(46, 4).Value2 = 45956.999976851854

$ws.Cells.Item(47, 1).Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(47, 2).Value2 = "003B号直流"
$ws.Cells.Item(47, 3).Value2 = 45956.444328703707
$ws.Cells.Item(47, 4).Value2 = 45956.999976851854

$ws.Cells.Item(48, 1).Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(48, 2).Value2 = "001A号直流"
$ws.Cells.Item(48, 3).Value2 = 45956.465254629627
$ws.Cells.Item(48, 4).Value2 = 45956.999976851854

$ws.Cells.Item(49, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(49, 2).Value2 = "B05号直流"
$ws.Cells.Item(49, 3).Value2 = 45956.47960648148
$ws.Cells.Item(49, 4).Value2 = 45956.999976851854

$ws.Cells.Item(50, 1).Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(50, 2).Value2 = "106号直流"
$ws.Cells.Item(50, 3).Value2 = 45956.488622685189
$ws.Cells.Item(50, 4).Value2 = 45956.999976851854

$ws.Cells.Item(51, 1).Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(51, 2).Value2 = "504号直流"
$ws.Cells.Item(51, 3).Value2 = 45956.491585648146
$ws.Cells.Item(51, 4).Value2 = 45956.999976851854

$ws.Cells.Item(52, 1).Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(52, 2).Value2 = "206号直流"
$ws.Cells.Item(52, 3).Value2 = 45956.496805555558
$ws.Cells.Item(52, 4).Value2 = 45956.999976851854

# Update sheet view: scroll back to top-left and set new selection to E20
$ws.Range("E20").Select()
